$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.926.92"
Set-TextValue $ws.Range("E2") "  +1.85%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.089.77"
Set-TextValue $ws.Range("E3") "  +1.18%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.07%  "

# Row 5
Set-TextValue $ws.Range("D5") "543.44"
Set-TextValue $ws.Range("E5") "  -0.28%  "

# Row 6
Set-TextValue $ws.Range("D6") "139.85"
Set-TextValue $ws.Range("E6") "  +4.37%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.14%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.083.69"
Set-TextValue $ws.Range("E8") "  +1.22%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.498"
Set-TextValue $ws.Range("E9") "  +2.00%  "

# Row 10
Set-TextValue $ws.Range("D10") "6.58"
Set-TextValue $ws.Range("E10") "  +3.15%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.156"
Set-TextValue $ws.Range("E11") "  +1.13%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.458"
Set-TextValue $ws.Range("E12") "  -0.06%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000226"
Set-TextValue $ws.Range("E13") "  +5.77%  "

# Row 14
Set-TextValue $ws.Range("D14") "34.75"
Set-TextValue $ws.Range("E14") "  +0.23%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.599.04"
Set-TextValue $ws.Range("E15") "  +1.55%  "

# Row 16
Set-TextValue $ws.Range("D16") "64.044.14"
Set-TextValue $ws.Range("E16") "  +1.93%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.112"
Set-TextValue $ws.Range("E17") "  +1.60%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.101.32"
Set-TextValue $ws.Range("E18") "  +1.57%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.66"
Set-TextValue $ws.Range("E19") "  +0.73%  "

# Row 20
Set-TextValue $ws.Range("D20") "481.30"
Set-TextValue $ws.Range("E20") "  -0.07%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.41"
Set-TextValue $ws.Range("E21") "  +0.66%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.700"
Set-TextValue $ws.Range("E22") "  +0.44%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.11"
Set-TextValue $ws.Range("E23") "  +1.05%  "

# Row 24
Set-TextValue $ws.Range("D24") "79.15"
Set-TextValue $ws.Range("E24") "  +2.63%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.39"
Set-TextValue $ws.Range("E25") "  +2.06%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.17%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.73"
Set-TextValue $ws.Range("E27") "  +1.17%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.07"
Set-TextValue $ws.Range("E28") "  -2.37%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  +0.15%  "

# Row 30
Set-TextValue $ws.Range("D30") "26.29"
Set-TextValue $ws.Range("E30") "  +0.81%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.90"
Set-TextValue $ws.Range("E31") "  -1.33%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +2.93%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.37"
Set-TextValue $ws.Range("E33") "  -4.71%  "

# Row 34
Set-TextValue $ws.Range("D34") "57.29"
Set-TextValue $ws.Range("E34") "  -1.78%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.39"
Set-TextValue $ws.Range("E35") "  +6.19%  "

# Row 36
Set-TextValue $ws.Range("D36") "495.42"
Set-TextValue $ws.Range("E36") "  -3.15%  "

# Row 37
Set-TextValue $ws.Range("D37") "6.01"
Set-TextValue $ws.Range("E37") "  +1.18%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.279.66"
Set-TextValue $ws.Range("E38") "  +6.65%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0405"
Set-TextValue $ws.Range("E39") "  +1.84%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0799"
Set-TextValue $ws.Range("E40") "  +1.70%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.120"
Set-TextValue $ws.Range("E41") "  +1.68%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.71"
Set-TextValue $ws.Range("E42") "  +4.57%  "

# Row 43
Set-TextValue $ws.Range("D43") "8.10"
Set-TextValue $ws.Range("E43") "  +1.12%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.254"
Set-TextValue $ws.Range("E44") "  +1.15%  "

# Row 46
Set-TextValue $ws.Range("B46") "Monero"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D46") "123.86"
Set-TextValue $ws.Range("E46") "  +3.19%  "

# Row 47
Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "25.22"
Set-TextValue $ws.Range("E47") "  +3.50%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.05"
Set-TextValue $ws.Range("E48") "  +1.01%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0₃0532"
Set-TextValue $ws.Range("E49") "  +8.28%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.109"
Set-TextValue $ws.Range("E50") "  +2.30%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.39"
Set-TextValue $ws.Range("E51") "  +0.14%  "

